$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(616).Insert()
$ws.Cells.Item(616, 1).Value = 3
$ws.Cells.Item(616, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(616, 3).Value = "Coquimbo"
$ws.Cells.Item(616, 4).Value = 45142
$ws.Cells.Item(616, 5).Value = 5
$ws.Cells.Item(616, 6).Value = 100112017
$ws.Cells.Item(616, 7).Value = "Apio"
$ws.Cells.Item(616, 8).Value = "Americana (o)"
$ws.Cells.Item(616, 9).Value = "Primera"
$ws.Cells.Item(616, 10).Value = 120
$ws.Cells.Item(616, 11).Value = 8000
$ws.Cells.Item(616, 12).Value = 8000
$ws.Cells.Item(616, 13).Value = 8000
$ws.Cells.Item(616, 14).Value = "$/docena de matas"
$ws.Cells.Item(616, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(616, 16).Value = 1333
$ws.Cells.Item(616, 17).Value = 6
$ws.Cells.Item(616, 18).Value = "Hortaliza"
